$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.558.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "'3.796.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'419.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'128.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.38%  "
$ws.Range("D7").Value = "'3.796.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.25%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.717"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.29%  "
$ws.Range("E11").Value = "  -7.90%  "
$ws.Range("D12").Value = "'0.0000342"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.56%  "
$ws.Range("E13").Value = "  -10.58%  "
$ws.Range("D14").Value = "'4.396.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "'15.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +19.48%  "
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "'3.775.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "'19.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.04%  "
$ws.Range("D20").Value = "'66.651.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  -4.53%  "
$ws.Range("D22").Value = "'402.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.56%  "
$ws.Range("D23").Value = "'14.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.34%  "
$ws.Range("D24").Value = "'83.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.50%  "
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("D26").Value = "'36.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("E27").Value = "  +10.55%  "
$ws.Range("D28").Value = "'3.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").Value = "'9.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.24%  "
$ws.Range("D30").Value = "'699.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "'8.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.27%  "
$ws.Range("D32").Value = "'2.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  -6.79%  "
$ws.Range("D37").Value = "'38.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.60%  "
$ws.Range("D38").Value = "'54.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.28%  "
$ws.Range("D39").Value = "'0.0₃0766"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.16%  "
$ws.Range("E40").Value = "  -6.60%  "
$ws.Range("D41").Value = "'2.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("D42").Value = "'4.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.29%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "'0.134"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.84%  "
$ws.Range("D45").Value = "'3.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("D46").Value = "'144.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("D47").Value = "'3.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.12%  "
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("D49").Value = "'25.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.11%  "
$ws.Range("D50").Value = "'2.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "'2.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.54%  "
